# scheduleGen: finish implementing the generated schedule — fix a couple of
# swapped subjects in the morning blocks, add the early-afternoon block to
# the first "turma" (352), and append three more full day-blocks (turmas
# 351 / 353 / extra) covering the afternoon shift, rows 31-60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fill-Row($row, $time, $b, $c, $d) {
    $ws.Cells.Item($row, 1).Value = $time
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
}

# Fills a standard 5-slot block (morning or afternoon) starting at $startRow
# with "---(---)" placeholders in every subject column, returns the next
# free row.
function Fill-EmptyBlock($startRow, $times) {
    $r = $startRow
    foreach ($t in $times) {
        Fill-Row $r $t "---(---)" "---(---)" "---(---)"
        $r = $r + 1
    }
    return $r
}

$amTimes = @("7h45", "8h30", "9h15", "10h15", "11h00")
$pmTimes = @("13h15", "14h00", "14h45", "15h45", "16h30")

# --- Turma 352 (col B) morning block: swap geogr/histo on a couple of rows ---
$ws.Cells.Item(2, 2).Value = "histo(marcio)"
$ws.Cells.Item(4, 2).Value = "geogr(marcio)"
$ws.Cells.Item(6, 2).Value = "geogr(marcio)"

# --- Row 7 separator becomes the lunch "INTERVALO" marker ---
$ws.Cells.Item(7, 1).Value = "INTERVALO"

# --- New early-afternoon block for turma 352, rows 8-12 ---
# (B10/C10/D10/C9/D9/C12/D12/C11/D11/D8 already hold the right placeholder
# values from the original template, so only touch what the diff changes.)
Fill-Row 8  "13h15" "histo(marcio)" "histo(evando)" "---(---)"
$ws.Cells.Item(9, 1).Value = "14h00"
$ws.Cells.Item(9, 2).Value = "geogr(marcio)"
$ws.Cells.Item(10, 1).Value = "14h45"
$ws.Cells.Item(11, 1).Value = "15h45"
$ws.Cells.Item(11, 2).Value = "histo(marcio)"
$ws.Cells.Item(12, 1).Value = "16h30"
$ws.Cells.Item(12, 2).Value = "histo(marcio)"

# --- Row 14: fill in the previously-empty subject for turma 352 ---
$ws.Cells.Item(14, 2).Value = "histo(marcio)"

# --- Row 19 separator becomes "INTERVALO" and times 20-24 shift to the PM slots ---
$ws.Cells.Item(19, 1).Value = "INTERVALO"
$ws.Cells.Item(20, 1).Value = "13h15"
$ws.Cells.Item(21, 1).Value = "14h00"
$ws.Cells.Item(22, 1).Value = "14h45"
$ws.Cells.Item(23, 1).Value = "15h45"
$ws.Cells.Item(24, 1).Value = "16h30"

# --- New afternoon-shift day blocks, rows 31-60 ---
$r = 31
$ws.Cells.Item($r, 1).Value = "INTERVALO"
$r = $r + 1
$r = Fill-EmptyBlock $r $pmTimes      # 32-36
$r = $r + 1                           # 37 blank separator
$r = Fill-EmptyBlock $r $amTimes      # 38-42
$ws.Cells.Item($r, 1).Value = "INTERVALO"   # 43
$r = $r + 1
$r = Fill-EmptyBlock $r $pmTimes      # 44-48
$r = $r + 1                           # 49 blank separator
$r = Fill-EmptyBlock $r $amTimes      # 50-54
$ws.Cells.Item($r, 1).Value = "INTERVALO"   # 55
$r = $r + 1
$r = Fill-EmptyBlock $r $pmTimes      # 56-60
